# completed bill view page
# - Add "Designation" header label in A1
# - Move the saved cell selection to F6
# - Switch the sheet to portrait page orientation for printing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the designation column
$ws.Range("A1").Value = "Designation"

# Set the page to print in portrait orientation
$ws.PageSetup.Orientation = 1

# Leave the saved selection on F6, matching where editing left off
$ws.Range("F6").Select() | Out-Null
